# Added periodic & upfront related scenarios
# Updates the "repaymentstrategy" answer on the ProductLoanInput sheet from
# "Mifos style" to "Penalties, Fees, Interest, Principal order" and gives the
# cell a left/top aligned style, then leaves the selection on that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

$cell.Select() | Out-Null
